$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.020676
$ws.Range("H2").Value = 0.062028
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2189473333333334
$ws.Range("N2").Value = 0.656842
$ws.Range("O2").Value = 0.009402596261870986
$ws.Range("P2").Value = 0.009402596261870984
$ws.Range("Q2").Value = 0.004526955064000001
$ws.Range("R2").Value = 0.040742595576
$ws.Range("S2").Value = 0.009402596261870986
$ws.Range("T2").Value = 0.009402596261870984

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.020676
$ws.Range("H3").Value = 0.062028
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.07911633333334
$ws.Range("N3").Value = 60.237349
$ws.Range("O3").Value = 0.8622887582286424
$ws.Range("P3").Value = 0.8622887582286423
$ws.Range("Q3").Value = 0.4151558093080001
$ws.Range("R3").Value = 3.736402283772
$ws.Range("S3").Value = 0.8622887582286424
$ws.Range("T3").Value = 0.8622887582286423

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.020676
$ws.Range("H4").Value = 0.062028
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.823530666666667
$ws.Range("N4").Value = 8.470592
$ws.Range("O4").Value = 0.1212552739852724
$ws.Range("P4").Value = 0.1212552739852723
$ws.Range("Q4").Value = 0.058379320064
$ws.Range("R4").Value = 0.525413880576
$ws.Range("S4").Value = 0.1212552739852724
$ws.Range("T4").Value = 0.1212552739852723

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.020676
$ws.Range("H5").Value = 0.062028
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1642436666666667
$ws.Range("N5").Value = 0.492731
$ws.Range("O5").Value = 0.007053371524214274
$ws.Range("P5").Value = 0.007053371524214274
$ws.Range("Q5").Value = 0.003395902052
$ws.Range("R5").Value = 0.030563118468
$ws.Range("S5").Value = 0.007053371524214274
$ws.Range("T5").Value = 0.007053371524214274
